# Apply the "complete project" data edits to the loginpage sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loginpage")

# Row 3: "renjini" -> "admin" ; numeric 12345 -> "xyzabcjkl"
$ws.Range("A3").Value2 = "admin"
$ws.Range("B3").Value2 = "xyzabcjkl"

# Row 4: "tester" -> "test1" ; "test123" -> "admin"
$ws.Range("A4").Value2 = "test1"
$ws.Range("B4").Value2 = "admin"

# Row 5: new row "sree" / "tester"
$ws.Range("A5").Value2 = "sree"
$ws.Range("B5").Value2 = "tester"

# Update the selected cell to match the final saved selection
$ws.Range("B3").Select() | Out-Null
